# Trade #104 closed at 2026-02-16 21:40:50 - leadlag UP +0.000%
#
# Two leadlag trades (rows 55/56 in the "leadlag" sheet, now also appended
# to "All Trades") and one momentum trade (row 18 in the "momentum" sheet,
# also appended to "All Trades") transition from OPEN -> CLOSED with their
# exit data filled in. A brand-new OPEN leadlag trade (#104) is appended to
# the "leadlag" sheet. The Summary and Comparison roll-up tables are
# refreshed to reflect the new totals.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1) "leadlag" sheet — close out trade #71 (row 55) and #72 (row 56)
# ---------------------------------------------------------------------
$wsLead = $wb.Worksheets.Item("leadlag")

# Row 55 -> trade #71
$wsLead.Cells.Item(55, 7).Value  = 69721.817427      # Exit Price (G)
$wsLead.Cells.Item(55, 8).Value  = "CLOSED"          # Status (H)
$wsLead.Cells.Item(55, 9).Value  = 1.6882            # P&L % (I)
$wsLead.Cells.Item(55, 10).Value = 16.88             # P&L $ (J)
$wsLead.Cells.Item(55, 13).Value = "time_exit_5min"  # Exit Reason (M)
$wsLead.Cells.Item(55, 14).Value = 5                 # Duration (N)

# Row 56 -> trade #72
$wsLead.Cells.Item(56, 7).Value  = 68496.225597
$wsLead.Cells.Item(56, 8).Value  = "CLOSED"
$wsLead.Cells.Item(56, 9).Value  = -0.1046
$wsLead.Cells.Item(56, 10).Value = -1.05
$wsLead.Cells.Item(56, 13).Value = "time_exit_5min"
$wsLead.Cells.Item(56, 14).Value = 5

# New trade #104, appended as row 80 (still OPEN). Copy/PasteSpecial the
# prior row first so the text-typed Date column ("2026-02-16") is carried
# over verbatim instead of being re-parsed (and auto-converted to a date
# serial number) by a fresh .Value assignment.
$wsLead.Rows.Item(79).Copy()
$wsLead.Rows.Item(80).PasteSpecial()
$wsLead.Cells.Item(80, 1).Value  = 104
$wsLead.Cells.Item(80, 3).Value  = "21:40:50"
$wsLead.Cells.Item(80, 6).Value  = 68479.035
$wsLead.Cells.Item(80, 12).Value = "Binance leading with 0.104% move"

# ---------------------------------------------------------------------
# 2) "momentum" sheet — close out trade #70 (row 18)
# ---------------------------------------------------------------------
$wsMom = $wb.Worksheets.Item("momentum")

$wsMom.Cells.Item(18, 7).Value  = 68068.900381
$wsMom.Cells.Item(18, 8).Value  = "CLOSED"
$wsMom.Cells.Item(18, 9).Value  = 0.656
$wsMom.Cells.Item(18, 10).Value = 6.56
$wsMom.Cells.Item(18, 13).Value = "time_exit_5min"
$wsMom.Cells.Item(18, 14).Value = 5

# ---------------------------------------------------------------------
# 3) "All Trades" sheet — append the three newly-closed trades
# ---------------------------------------------------------------------
$wsAll = $wb.Worksheets.Item("All Trades")

# Row 71 <- momentum trade #70
$wsAll.Rows.Item(70).Copy()
$wsAll.Rows.Item(71).PasteSpecial()
$wsAll.Cells.Item(71, 1).Value  = 70
$wsAll.Cells.Item(71, 3).Value  = "21:35:31"
$wsAll.Cells.Item(71, 6).Value  = 68518.41499999999
$wsAll.Cells.Item(71, 7).Value  = 68068.900381
$wsAll.Cells.Item(71, 9).Value  = 0.656
$wsAll.Cells.Item(71, 10).Value = 6.56
$wsAll.Cells.Item(71, 12).Value = "Downward momentum: -0.291% over 10 samples"

# Row 72 <- leadlag trade #71
$wsAll.Rows.Item(70).Copy()
$wsAll.Rows.Item(72).PasteSpecial()
$wsAll.Cells.Item(72, 1).Value  = 71
$wsAll.Cells.Item(72, 3).Value  = "21:35:38"
$wsAll.Cells.Item(72, 4).Value  = "leadlag"
$wsAll.Cells.Item(72, 5).Value  = "UP"
$wsAll.Cells.Item(72, 6).Value  = 68564.31
$wsAll.Cells.Item(72, 7).Value  = 69721.817427
$wsAll.Cells.Item(72, 9).Value  = 1.6882
$wsAll.Cells.Item(72, 10).Value = 16.88
$wsAll.Cells.Item(72, 11).Value = 0.75
$wsAll.Cells.Item(72, 12).Value = "Binance leading with 0.084% move"

# Row 73 <- leadlag trade #72
$wsAll.Rows.Item(70).Copy()
$wsAll.Rows.Item(73).PasteSpecial()
$wsAll.Cells.Item(73, 1).Value  = 72
$wsAll.Cells.Item(73, 3).Value  = "21:35:44"
$wsAll.Cells.Item(73, 4).Value  = "leadlag"
$wsAll.Cells.Item(73, 5).Value  = "UP"
$wsAll.Cells.Item(73, 6).Value  = 68567.925
$wsAll.Cells.Item(73, 7).Value  = 68496.225597
$wsAll.Cells.Item(73, 9).Value  = -0.1046
$wsAll.Cells.Item(73, 10).Value = -1.05
$wsAll.Cells.Item(73, 11).Value = 0.75
$wsAll.Cells.Item(73, 12).Value = "Coinbase leading with 0.085% move"

# ---------------------------------------------------------------------
# 4) "Summary" sheet — refresh OVERALL / leadlag / momentum roll-ups
# ---------------------------------------------------------------------
$wsSum = $wb.Worksheets.Item("Summary")

$wsSum.Cells.Item(2, 3).Value = 72          # OVERALL Total Trades
$wsSum.Cells.Item(2, 5).Value = "+21.4923%" # OVERALL Total P&L %
$wsSum.Cells.Item(2, 6).Value = "+0.2985%"  # OVERALL Avg Trade

$wsSum.Cells.Item(3, 3).Value = 78          # leadlag Total Trades
$wsSum.Cells.Item(3, 4).Value = "43.6%"     # leadlag Win Rate
$wsSum.Cells.Item(3, 5).Value = "+12.6749%" # leadlag Total P&L %
$wsSum.Cells.Item(3, 6).Value = "+0.1625%"  # leadlag Avg Trade

$wsSum.Cells.Item(4, 4).Value = "60.0%"     # momentum Win Rate
$wsSum.Cells.Item(4, 5).Value = "+8.8174%"  # momentum Total P&L %
$wsSum.Cells.Item(4, 6).Value = "+0.3527%"  # momentum Avg Trade

# ---------------------------------------------------------------------
# 5) "Comparison" sheet — refresh leadlag / momentum stats
# ---------------------------------------------------------------------
$wsCmp = $wb.Worksheets.Item("Comparison")

$wsCmp.Cells.Item(2, 2).Value = 78          # leadlag Total Trades
$wsCmp.Cells.Item(2, 3).Value = "43.6%"     # leadlag Win Rate
$wsCmp.Cells.Item(2, 4).Value = "3.06"      # leadlag Profit Factor
$wsCmp.Cells.Item(2, 5).Value = "+0.5540%"  # leadlag Avg Win %
$wsCmp.Cells.Item(2, 6).Value = "-0.2933%"  # leadlag Avg Loss %
$wsCmp.Cells.Item(2, 7).Value = "1.89"      # leadlag Win/Loss Ratio

$wsCmp.Cells.Item(3, 3).Value = "60.0%"     # momentum Win Rate
$wsCmp.Cells.Item(3, 4).Value = "8.84"      # momentum Profit Factor
$wsCmp.Cells.Item(3, 5).Value = "+0.6628%"  # momentum Avg Win %

Write-Output "edit applied"
